$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.209.61'
$ws.Cells.Item(2, 5).Value = '  +1.37%  '
$ws.Cells.Item(3, 4).Value = '1.808.45'
$ws.Cells.Item(3, 5).Value = '  +2.31%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  -0.17%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '338.22'
$ws.Cells.Item(5, 5).Value = '  -0.04%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9973'
$ws.Cells.Item(6, 5).Value = '  -0.38%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4579'
$ws.Cells.Item(7, 5).Value = '  +21.12%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3557'
$ws.Cells.Item(8, 5).Value = '  +5.27%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '45.61'
$ws.Cells.Item(9, 5).Value = '  +0.25%  '
$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.07661'
$ws.Cells.Item(10, 5).Value = '  +5.00%  '
$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.155'
$ws.Cells.Item(11, 5).Value = '  +1.96%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '22.87'
$ws.Cells.Item(12, 5).Value = '  -1.36%  '
$ws.Cells.Item(13, 5).Value = '  -0.08%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.280'
$ws.Cells.Item(14, 5).Value = '  +0.01%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.282'
$ws.Cells.Item(15, 5).Value = '  -0.08%  '
$ws.Cells.Item(16, 4).Value = '1.806.40'
$ws.Cells.Item(16, 5).Value = '  +2.24%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.00001102'
$ws.Cells.Item(17, 5).Value = '  +4.22%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.06691'
$ws.Cells.Item(18, 5).Value = '  +1.26%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '82.09'
$ws.Cells.Item(19, 5).Value = '  +1.11%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.9989'
$ws.Cells.Item(20, 5).Value = '  -0.31%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '17.30'
$ws.Cells.Item(21, 5).Value = '  +0.46%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.412'
$ws.Cells.Item(22, 5).Value = '  +0.85%  '
$ws.Cells.Item(23, 4).Value = '28.255.62'
$ws.Cells.Item(23, 5).Value = '  +1.46%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '12.00'
$ws.Cells.Item(24, 5).Value = '  +1.32%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.405'
$ws.Cells.Item(25, 5).Value = '  +0.87%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '20.85'
$ws.Cells.Item(26, 5).Value = '  +3.46%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.437'
$ws.Cells.Item(27, 5).Value = '  +3.24%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '155.76'
$ws.Cells.Item(28, 5).Value = '  +2.80%  '
$ws.Cells.Item(29, 4).Value = '2.012.24'
$ws.Cells.Item(29, 5).Value = '  +2.23%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.302'
$ws.Cells.Item(30, 5).Value = '  -14.45%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '134.11'
$ws.Cells.Item(31, 5).Value = '  +0.63%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.070'
$ws.Cells.Item(32, 5).Value = '  +0.79%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.979'
$ws.Cells.Item(33, 5).Value = '  +0.56%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.09514'
$ws.Cells.Item(34, 5).Value = '  +8.44%  '
$ws.Cells.Item(35, 2).Value = 'VeChain'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.02390'
$ws.Cells.Item(35, 5).Value = '  +0.77%  '
$ws.Cells.Item(36, 2).Value = 'Aptos'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '12.27'
$ws.Cells.Item(36, 5).Value = '  -1.46%  '
$ws.Cells.Item(37, 2).Value = 'TheSandbox'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.6804'
$ws.Cells.Item(37, 5).Value = '  +1.00%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.06274'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.2177'
$ws.Cells.Item(39, 5).Value = '  +2.38%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.217'
$ws.Cells.Item(40, 5).Value = '  +0.10%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.491'
$ws.Cells.Item(41, 5).Value = '  +0.98%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.220'
$ws.Cells.Item(42, 5).Value = '  -0.21%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '8.194'
$ws.Cells.Item(43, 5).Value = '  +1.27%  '
$ws.Cells.Item(44, 2).Value = 'Frax'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.9975'
$ws.Cells.Item(44, 5).Value = '  -0.35%  '
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '14.05'
$ws.Cells.Item(45, 5).Value = '  +1.31%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.6153'
$ws.Cells.Item(46, 5).Value = '  +0.10%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.865'
$ws.Cells.Item(47, 5).Value = '  +0.47%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '129.93'
$ws.Cells.Item(48, 5).Value = '  -1.44%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.048'
$ws.Cells.Item(49, 5).Value = '  +0.55%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.07122'
$ws.Cells.Item(50, 5).Value = '  -2.16%  '
$ws.Cells.Item(51, 5).Value = '  -1.80%  '
